$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.0244
$ws.Range("E2").Value = 0.9756

$ws.Range("D3").Value = 0.8077
$ws.Range("E3").Value = 0.1923

$ws.Range("D4").Value = 0.0244
$ws.Range("E4").Value = 0.9756

$ws.Range("D5").Value = 0.0244
$ws.Range("E5").Value = 0.9756

$ws.Range("D6").Value = 0.0193
$ws.Range("E6").Value = 0.9807

$ws.Range("D7").Value = 0.0514
$ws.Range("E7").Value = 0.9486

$ws.Range("D8").Value = 0.0197
$ws.Range("E8").Value = 0.9803

$ws.Range("D9").Value = 0.0244
$ws.Range("E9").Value = 0.9756

$ws.Range("D10").Value = 0.0244
$ws.Range("E10").Value = 0.9756

$ws.Range("D11").Value = 0.0244
$ws.Range("E11").Value = 0.9756
